# Update the 5x20 "addition/subtraction within 100" answer table: every
# cell's arithmetic expression is replaced with a new one. Every old
# expression is unique in the document, so a plain Find/Replace (no
# wildcards, match whole string) targets exactly one cell each time.
#
# One pair of cells forms a short dependency chain (row 9, cols 2 & 4):
# the text "36+30=66" is both being replaced *and* being written as a
# *new* value elsewhere, so the replacement that turns the original
# "36+30=66" cell into "37+47=84" must run before the replacement that
# writes a fresh "36+30=66" into another cell (otherwise the later
# Find could match the freshly written text instead of the original).
# The order below already accounts for that.
$d = $word.ActiveDocument

$d.Content.Find.Execute("43+29=72", $true, $false, $false, $false, $false, $true, 1, $false, "68-25=43", 2) | Out-Null
$d.Content.Find.Execute("55-26=29", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=6", 2) | Out-Null
$d.Content.Find.Execute("64+7=71", $true, $false, $false, $false, $false, $true, 1, $false, "29+47=76", 2) | Out-Null
$d.Content.Find.Execute("83-69=14", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "72-35=37", 2) | Out-Null
$d.Content.Find.Execute("54+31=85", $true, $false, $false, $false, $false, $true, 1, $false, "51+45=96", 2) | Out-Null
$d.Content.Find.Execute("98-85=13", $true, $false, $false, $false, $false, $true, 1, $false, "14+2=16", 2) | Out-Null
$d.Content.Find.Execute("97-18=79", $true, $false, $false, $false, $false, $true, 1, $false, "87+1=88", 2) | Out-Null
$d.Content.Find.Execute("83+2=85", $true, $false, $false, $false, $false, $true, 1, $false, "39+6=45", 2) | Out-Null
$d.Content.Find.Execute("83-59=24", $true, $false, $false, $false, $false, $true, 1, $false, "89-52=37", 2) | Out-Null
$d.Content.Find.Execute("56+21=77", $true, $false, $false, $false, $false, $true, 1, $false, "2+1=3", 2) | Out-Null
$d.Content.Find.Execute("81+13=94", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=63", 2) | Out-Null
$d.Content.Find.Execute("68-12=56", $true, $false, $false, $false, $false, $true, 1, $false, "59+23=82", 2) | Out-Null
$d.Content.Find.Execute("77+22=99", $true, $false, $false, $false, $false, $true, 1, $false, "79-30=49", 2) | Out-Null
$d.Content.Find.Execute("60-59=1", $true, $false, $false, $false, $false, $true, 1, $false, "88-60=28", 2) | Out-Null
$d.Content.Find.Execute("53-5=48", $true, $false, $false, $false, $false, $true, 1, $false, "52-39=13", 2) | Out-Null
$d.Content.Find.Execute("75-64=11", $true, $false, $false, $false, $false, $true, 1, $false, "60+28=88", 2) | Out-Null
$d.Content.Find.Execute("71-41=30", $true, $false, $false, $false, $false, $true, 1, $false, "63-9=54", 2) | Out-Null
$d.Content.Find.Execute("42+22=64", $true, $false, $false, $false, $false, $true, 1, $false, "62-35=27", 2) | Out-Null
$d.Content.Find.Execute("46+1=47", $true, $false, $false, $false, $false, $true, 1, $false, "35+62=97", 2) | Out-Null
$d.Content.Find.Execute("29+63=92", $true, $false, $false, $false, $false, $true, 1, $false, "31+49=80", 2) | Out-Null
$d.Content.Find.Execute("99-5=94", $true, $false, $false, $false, $false, $true, 1, $false, "6+46=52", 2) | Out-Null
$d.Content.Find.Execute("28+63=91", $true, $false, $false, $false, $false, $true, 1, $false, "66-60=6", 2) | Out-Null
$d.Content.Find.Execute("40+27=67", $true, $false, $false, $false, $false, $true, 1, $false, "68-20=48", 2) | Out-Null
$d.Content.Find.Execute("96-46=50", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=7", 2) | Out-Null
$d.Content.Find.Execute("85-54=31", $true, $false, $false, $false, $false, $true, 1, $false, "74-63=11", 2) | Out-Null
$d.Content.Find.Execute("34-1=33", $true, $false, $false, $false, $false, $true, 1, $false, "84-0=84", 2) | Out-Null
$d.Content.Find.Execute("20+78=98", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("3+0=3", $true, $false, $false, $false, $false, $true, 1, $false, "62+5=67", 2) | Out-Null
$d.Content.Find.Execute("7+70=77", $true, $false, $false, $false, $false, $true, 1, $false, "96-60=36", 2) | Out-Null
$d.Content.Find.Execute("47-10=37", $true, $false, $false, $false, $false, $true, 1, $false, "31+10=41", 2) | Out-Null
$d.Content.Find.Execute("97-58=39", $true, $false, $false, $false, $false, $true, 1, $false, "48-4=44", 2) | Out-Null
$d.Content.Find.Execute("68+6=74", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=9", 2) | Out-Null
$d.Content.Find.Execute("57-4=53", $true, $false, $false, $false, $false, $true, 1, $false, "29+38=67", 2) | Out-Null
$d.Content.Find.Execute("34+33=67", $true, $false, $false, $false, $false, $true, 1, $false, "82-52=30", 2) | Out-Null
$d.Content.Find.Execute("65+3=68", $true, $false, $false, $false, $false, $true, 1, $false, "24-3=21", 2) | Out-Null
$d.Content.Find.Execute("82+16=98", $true, $false, $false, $false, $false, $true, 1, $false, "65+14=79", 2) | Out-Null
$d.Content.Find.Execute("82-50=32", $true, $false, $false, $false, $false, $true, 1, $false, "58-1=57", 2) | Out-Null
$d.Content.Find.Execute("11+45=56", $true, $false, $false, $false, $false, $true, 1, $false, "89+8=97", 2) | Out-Null
$d.Content.Find.Execute("52-17=35", $true, $false, $false, $false, $false, $true, 1, $false, "93-62=31", 2) | Out-Null
$d.Content.Find.Execute("14+33=47", $true, $false, $false, $false, $false, $true, 1, $false, "91+7=98", 2) | Out-Null
$d.Content.Find.Execute("74-41=33", $true, $false, $false, $false, $false, $true, 1, $false, "66+33=99", 2) | Out-Null
$d.Content.Find.Execute("36+30=66", $true, $false, $false, $false, $false, $true, 1, $false, "37+47=84", 2) | Out-Null
$d.Content.Find.Execute("13+41=54", $true, $false, $false, $false, $false, $true, 1, $false, "36+30=66", 2) | Out-Null
$d.Content.Find.Execute("24+68=92", $true, $false, $false, $false, $false, $true, 1, $false, "92-47=45", 2) | Out-Null
$d.Content.Find.Execute("10+34=44", $true, $false, $false, $false, $false, $true, 1, $false, "89-85=4", 2) | Out-Null
$d.Content.Find.Execute("29+67=96", $true, $false, $false, $false, $false, $true, 1, $false, "82-13=69", 2) | Out-Null
$d.Content.Find.Execute("54-48=6", $true, $false, $false, $false, $false, $true, 1, $false, "90-87=3", 2) | Out-Null
$d.Content.Find.Execute("20+1=21", $true, $false, $false, $false, $false, $true, 1, $false, "7+37=44", 2) | Out-Null
$d.Content.Find.Execute("32+57=89", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=22", 2) | Out-Null
$d.Content.Find.Execute("41+9=50", $true, $false, $false, $false, $false, $true, 1, $false, "64-36=28", 2) | Out-Null
$d.Content.Find.Execute("45+54=99", $true, $false, $false, $false, $false, $true, 1, $false, "30+40=70", 2) | Out-Null
$d.Content.Find.Execute("60+14=74", $true, $false, $false, $false, $false, $true, 1, $false, "36+1=37", 2) | Out-Null
$d.Content.Find.Execute("84-58=26", $true, $false, $false, $false, $false, $true, 1, $false, "10+62=72", 2) | Out-Null
$d.Content.Find.Execute("75+17=92", $true, $false, $false, $false, $false, $true, 1, $false, "35+19=54", 2) | Out-Null
$d.Content.Find.Execute("71-63=8", $true, $false, $false, $false, $false, $true, 1, $false, "34+36=70", 2) | Out-Null
$d.Content.Find.Execute("86-12=74", $true, $false, $false, $false, $false, $true, 1, $false, "71-67=4", 2) | Out-Null
$d.Content.Find.Execute("25+19=44", $true, $false, $false, $false, $false, $true, 1, $false, "99-43=56", 2) | Out-Null
$d.Content.Find.Execute("31+28=59", $true, $false, $false, $false, $false, $true, 1, $false, "24+46=70", 2) | Out-Null
$d.Content.Find.Execute("6+89=95", $true, $false, $false, $false, $false, $true, 1, $false, "77+1=78", 2) | Out-Null
$d.Content.Find.Execute("97-19=78", $true, $false, $false, $false, $false, $true, 1, $false, "81+8=89", 2) | Out-Null
$d.Content.Find.Execute("13+20=33", $true, $false, $false, $false, $false, $true, 1, $false, "33+19=52", 2) | Out-Null
$d.Content.Find.Execute("70-37=33", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("9+15=24", $true, $false, $false, $false, $false, $true, 1, $false, "80+3=83", 2) | Out-Null
$d.Content.Find.Execute("60-40=20", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=7", 2) | Out-Null
$d.Content.Find.Execute("67-26=41", $true, $false, $false, $false, $false, $true, 1, $false, "53+3=56", 2) | Out-Null
$d.Content.Find.Execute("67-49=18", $true, $false, $false, $false, $false, $true, 1, $false, "81-65=16", 2) | Out-Null
$d.Content.Find.Execute("29-17=12", $true, $false, $false, $false, $false, $true, 1, $false, "20+34=54", 2) | Out-Null
$d.Content.Find.Execute("70+7=77", $true, $false, $false, $false, $false, $true, 1, $false, "38-34=4", 2) | Out-Null
$d.Content.Find.Execute("64-32=32", $true, $false, $false, $false, $false, $true, 1, $false, "8+60=68", 2) | Out-Null
$d.Content.Find.Execute("72-29=43", $true, $false, $false, $false, $false, $true, 1, $false, "73-71=2", 2) | Out-Null
$d.Content.Find.Execute("38+6=44", $true, $false, $false, $false, $false, $true, 1, $false, "59-26=33", 2) | Out-Null
$d.Content.Find.Execute("85+3=88", $true, $false, $false, $false, $false, $true, 1, $false, "45+44=89", 2) | Out-Null
$d.Content.Find.Execute("62-53=9", $true, $false, $false, $false, $false, $true, 1, $false, "75-61=14", 2) | Out-Null
$d.Content.Find.Execute("18+54=72", $true, $false, $false, $false, $false, $true, 1, $false, "80-60=20", 2) | Out-Null
$d.Content.Find.Execute("49+47=96", $true, $false, $false, $false, $false, $true, 1, $false, "64-7=57", 2) | Out-Null
$d.Content.Find.Execute("25+53=78", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=55", 2) | Out-Null
$d.Content.Find.Execute("43-40=3", $true, $false, $false, $false, $false, $true, 1, $false, "13+46=59", 2) | Out-Null
$d.Content.Find.Execute("77-19=58", $true, $false, $false, $false, $false, $true, 1, $false, "19+60=79", 2) | Out-Null
$d.Content.Find.Execute("78-12=66", $true, $false, $false, $false, $false, $true, 1, $false, "74-50=24", 2) | Out-Null
$d.Content.Find.Execute("45-39=6", $true, $false, $false, $false, $false, $true, 1, $false, "70-14=56", 2) | Out-Null
$d.Content.Find.Execute("99-2=97", $true, $false, $false, $false, $false, $true, 1, $false, "3+7=10", 2) | Out-Null
$d.Content.Find.Execute("13+64=77", $true, $false, $false, $false, $false, $true, 1, $false, "92-50=42", 2) | Out-Null
$d.Content.Find.Execute("12+83=95", $true, $false, $false, $false, $false, $true, 1, $false, "24+55=79", 2) | Out-Null
$d.Content.Find.Execute("59-12=47", $true, $false, $false, $false, $false, $true, 1, $false, "4+51=55", 2) | Out-Null
$d.Content.Find.Execute("25-3=22", $true, $false, $false, $false, $false, $true, 1, $false, "44+8=52", 2) | Out-Null
$d.Content.Find.Execute("58+11=69", $true, $false, $false, $false, $false, $true, 1, $false, "51+30=81", 2) | Out-Null
$d.Content.Find.Execute("80-54=26", $true, $false, $false, $false, $false, $true, 1, $false, "13-9=4", 2) | Out-Null
$d.Content.Find.Execute("46+52=98", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=77", 2) | Out-Null
$d.Content.Find.Execute("95-75=20", $true, $false, $false, $false, $false, $true, 1, $false, "76-8=68", 2) | Out-Null
$d.Content.Find.Execute("87-3=84", $true, $false, $false, $false, $false, $true, 1, $false, "12+46=58", 2) | Out-Null
$d.Content.Find.Execute("0+3=3", $true, $false, $false, $false, $false, $true, 1, $false, "52-29=23", 2) | Out-Null
$d.Content.Find.Execute("81-46=35", $true, $false, $false, $false, $false, $true, 1, $false, "30+26=56", 2) | Out-Null
$d.Content.Find.Execute("24+20=44", $true, $false, $false, $false, $false, $true, 1, $false, "87-53=34", 2) | Out-Null
$d.Content.Find.Execute("22-5=17", $true, $false, $false, $false, $false, $true, 1, $false, "93-92=1", 2) | Out-Null
$d.Content.Find.Execute("55+36=91", $true, $false, $false, $false, $false, $true, 1, $false, "15+24=39", 2) | Out-Null
$d.Content.Find.Execute("82-6=76", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=36", 2) | Out-Null
$d.Content.Find.Execute("82+4=86", $true, $false, $false, $false, $false, $true, 1, $false, "74+11=85", 2) | Out-Null
$d.Content.Find.Execute("39-32=7", $true, $false, $false, $false, $false, $true, 1, $false, "80-3=77", 2) | Out-Null
$d.Content.Find.Execute("51-24=27", $true, $false, $false, $false, $false, $true, 1, $false, "4+42=46", 2) | Out-Null

Write-Output "Replaced 100 cells"
